$wb = $excel.ActiveWorkbook

# --- Update the daily conversion note on sheet "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value2
$text = $text.Replace("✅ 1000 Bs = 1.28 = 4811.42 pesos", "✅ 1000 Bs = 1.38 = 5139.65 pesos")
$text = $text.Replace("✅ 4811.42 pesos = 1.29 = 901.39 Bs", "✅ 5139.65 pesos = 1.37 = 836.27 Bs")
$cell.Value2 = $text

# --- Update the rate figures on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 726.995
$ws2.Range("O10").Value = 3736.5
$ws2.Range("N12").Value = 3749
$ws2.Range("O12").Value = 610
